$wb = $excel.ActiveWorkbook

# --- "transitions" sheet (sheet1): update the treatment -> death transition
# from a "time_dependent" parameterisation to a "time_dependent_gompertz" one,
# with new parameter_1/parameter_2 values.
$wsTransitions = $wb.Worksheets.Item("transitions")
$wsTransitions.Range("C3").Value = "time_dependent_gompertz"
$wsTransitions.Range("D3").Value = -1.3624080000000001
$wsTransitions.Range("E3").Value = -0.33633940000000001

# --- "costs" sheet (sheet2): insert a new "type" column (B) describing how
# each state's cost is derived, shifting cost / cost_variance one column right.
$wsCosts = $wb.Worksheets.Item("costs")
$wsCosts.Columns("B:B").Insert()
$wsCosts.Range("B1").Value = "type"
$wsCosts.Range("B2").Value = "static"
$wsCosts.Range("B3").Value = "static"

# --- "utilities" sheet (sheet3): same "type" column insertion as costs.
$wsUtilities = $wb.Worksheets.Item("utilities")
$wsUtilities.Columns("B:B").Insert()
$wsUtilities.Range("B1").Value = "type"
$wsUtilities.Range("B2").Value = "static"
$wsUtilities.Range("B3").Value = "static"

# --- Restore per-sheet selections, then leave "transitions" as the active
# (visible) sheet/selection, matching the saved workbook view state.
$wsCosts.Activate()
$wsCosts.Range("B3").Select()

$wsUtilities.Activate()
$wsUtilities.Range("E5").Select()

$wsSpecification = $wb.Worksheets.Item("specification")
$wsSpecification.Activate()
$wsSpecification.Range("B13").Select()

$wsTransitions.Activate()
$wsTransitions.Range("G6").Select()
